$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3303.5
$ws.Range("J76").Value = 3291
$ws.Range("L76").Value = 3291
$ws.Range("N76").Value = -3921
$ws.Range("H79").Value = 3303.5
$ws.Range("J79").Value = 3291
$ws.Range("L79").Value = 3291
$ws.Range("N79").Value = -5475
$ws.Range("H88").Value = 3463.1333
$ws.Range("J88").Value = 3456.5
$ws.Range("L88").Value = 3456.5
$ws.Range("N88").Value = -4268.5
$ws.Range("H91").Value = 3463.1333
$ws.Range("J91").Value = 3456.5
$ws.Range("L91").Value = 3456.5
$ws.Range("N91").Value = -6264.5
$ws.Range("H138").Value = 2358.41
$ws.Range("I138").Value = 1308.9
$ws.Range("K138").Value = 3926.7
$ws.Range("M138").Value = 1213.3

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H86").Value = 59994.5
$ws.Range("J86").Value = 59994.5
$ws.Range("L86").Value = 59994.5
$ws.Range("N86").Value = -62366.5
$ws.Range("H89").Value = 59994.5
$ws.Range("J89").Value = 59994.5
$ws.Range("L89").Value = 179983.5
$ws.Range("N89").Value = -191839.5
$ws.Range("H97").Value = 1465.7
$ws.Range("I97").Value = 987
$ws.Range("J97").Value = 2782.125
$ws.Range("K97").Value = 987
$ws.Range("L97").Value = 2782.125
$ws.Range("M97").Value = -491
$ws.Range("N97").Value = -3774.125
$ws.Range("H132").Value = 3879.3794
$ws.Range("I132").Value = 2972.4424
$ws.Range("J132").Value = 11739.5
$ws.Range("K132").Value = 8917.3272
$ws.Range("L132").Value = 35218.5
$ws.Range("M132").Value = -6387.3272
$ws.Range("N132").Value = -40278.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 38743.5
$ws.Range("J2").Value = 37488
$ws.Range("L2").Value = 37488
$ws.Range("N2").Value = -37714
$ws.Range("H20").Value = 2890.32
$ws.Range("I20").Value = 1973.9445
$ws.Range("K20").Value = 1973.9445
$ws.Range("M20").Value = -1726.9445
$ws.Range("H86").Value = 1946.2142
$ws.Range("I86").Value = 1946.2142
$ws.Range("K86").Value = 1946.2142
$ws.Range("M86").Value = -823.2141999999999
$ws.Range("H89").Value = 1946.2142
$ws.Range("I89").Value = 1946.2142
$ws.Range("K89").Value = 9731.071
$ws.Range("M89").Value = -4115.071

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 484
$ws.Range("I22").Value = 398.5
$ws.Range("K22").Value = 398.5
$ws.Range("M22").Value = -48.5
$ws.Range("H58").Value = 1395.9565
$ws.Range("I58").Value = 1395.9565
$ws.Range("K58").Value = 1395.9565
$ws.Range("M58").Value = -1192.9565
$ws.Range("H99").Value = 1002891.6
$ws.Range("I99").Value = 1002891.6
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1002891.6
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -1001393.6
$ws.Range("N99").ClearContents()
$ws.Range("H126").Value = 1002891.6
$ws.Range("I126").Value = 1002891.6
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3008674.8
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3006204.8
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 1395.9565
$ws.Range("I136").Value = 1395.9565
$ws.Range("K136").Value = 4187.8695
$ws.Range("M136").Value = -1637.8695

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 1646.7778
$ws.Range("I39").Value = 826.7778
$ws.Range("J39").Value = 1920.1111
$ws.Range("K39").Value = 2480.3334
$ws.Range("L39").Value = 5760.3333
$ws.Range("M39").Value = -2186.3334
$ws.Range("N39").Value = -6348.3333
$ws.Range("H61").Value = 4246.1665
$ws.Range("I61").Value = 80
$ws.Range("J61").Value = 5079.4
$ws.Range("K61").Value = 240
$ws.Range("L61").Value = 15238.2
$ws.Range("M61").Value = -25
$ws.Range("N61").Value = -15668.2
$ws.Range("H121").Value = 2617.7693
$ws.Range("I121").Value = 1438.625
$ws.Range("J121").Value = 4504.4
$ws.Range("K121").Value = 4315.875
$ws.Range("L121").Value = 13513.2
$ws.Range("M121").Value = -3005.875
$ws.Range("N121").Value = -16133.2
$ws.Range("H132").Value = 6571.773
$ws.Range("I132").Value = 6571.773
$ws.Range("K132").Value = 59145.957
$ws.Range("M132").Value = -56615.957
$ws.Range("H137").Value = 2683
$ws.Range("I137").Value = 1766
$ws.Range("K137").Value = 5298
$ws.Range("M137").Value = -198

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 89.25
$ws.Range("I2").Value = 81.2
$ws.Range("J2").Value = 97.3
$ws.Range("K2").Value = 81.2
$ws.Range("L2").Value = 97.3
$ws.Range("M2").Value = 31.8
$ws.Range("N2").Value = -323.3
$ws.Range("H80").Value = 2814.5667
$ws.Range("I80").Value = 2697.4583
$ws.Range("K80").Value = 2697.4583
$ws.Range("M80").Value = -1699.4583
$ws.Range("H83").Value = 2814.5667
$ws.Range("I83").Value = 2697.4583
$ws.Range("K83").Value = 13487.2915
$ws.Range("M83").Value = -8495.291499999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6765.5454
$ws.Range("I7").Value = 6632.8423
$ws.Range("K7").Value = 6632.8423
$ws.Range("M7").Value = -6520.8423
$ws.Range("H40").Value = 17728.934
$ws.Range("I40").Value = 17728.934
$ws.Range("K40").Value = 17728.934
$ws.Range("M40").Value = -17592.934
$ws.Range("H61").Value = 1899.6
$ws.Range("I61").Value = 1499.6666
$ws.Range("J61").Value = 2499.5
$ws.Range("K61").Value = 1499.6666
$ws.Range("L61").Value = 2499.5
$ws.Range("M61").Value = -1297.6666
$ws.Range("N61").Value = -2903.5
$ws.Range("H100").Value = 11113377
$ws.Range("I100").Value = 14287700
$ws.Range("J100").Value = 3245
$ws.Range("K100").Value = 14287700
$ws.Range("L100").Value = 3245
$ws.Range("M100").Value = -14287159
$ws.Range("N100").Value = -4327
$ws.Range("H113").Value = 1899.6
$ws.Range("I113").Value = 1499.6666
$ws.Range("J113").Value = 2499.5
$ws.Range("K113").Value = 1499.6666
$ws.Range("L113").Value = 2499.5
$ws.Range("M113").Value = 670.3334
$ws.Range("N113").Value = -6839.5
$ws.Range("H126").Value = 6765.5454
$ws.Range("I126").Value = 6632.8423
$ws.Range("K126").Value = 19898.5269
$ws.Range("M126").Value = -17428.5269

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9012.846
$ws.Range("I62").Value = 4499.6665
$ws.Range("J62").Value = 10366.8
$ws.Range("K62").Value = 4499.6665
$ws.Range("L62").Value = 10366.8
$ws.Range("M62").Value = -3875.6665
$ws.Range("N62").Value = -11614.8
$ws.Range("H65").Value = 9012.846
$ws.Range("I65").Value = 4499.6665
$ws.Range("J65").Value = 10366.8
$ws.Range("K65").Value = 22498.3325
$ws.Range("L65").Value = 51834
$ws.Range("M65").Value = -19378.3325
$ws.Range("N65").Value = -58074
$ws.Range("H122").Value = 3674.611
$ws.Range("I122").Value = 3642.6428
$ws.Range("J122").Value = 3786.5
$ws.Range("K122").Value = 10927.9284
$ws.Range("L122").Value = 11359.5
$ws.Range("M122").Value = -8477.928400000001
$ws.Range("N122").Value = -16259.5